$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3283
$ws1.Range("F5").Value = 2364
$ws1.Range("F6").Value = 333
$ws1.Range("F7").Value = 329
$ws1.Range("F8").Value = 1339
$ws1.Range("F10").Value = 283
$ws1.Range("F11").Value = 496
$ws1.Range("F15").Value = 545
$ws1.Range("F16").Value = 8261
$ws1.Range("F17").Value = 359
$ws1.Range("F19").Value = 236
$ws1.Range("F23").Value = 563
$ws1.Range("F25").Value = 1147
$ws1.Range("F26").Value = 992
$ws1.Range("F27").Value = 1916
$ws1.Range("F28").Value = 1457
$ws1.Range("F30").Value = 238
$ws1.Range("F34").Value = 18
$ws1.Range("F37").Value = 297
$ws1.Range("F38").Value = 51
$ws1.Range("F39").Value = 213
$ws1.Range("F40").Value = 388
$ws1.Range("F41").Value = 44

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3283
$ws4.Range("F6").Value = 2364
$ws4.Range("F7").Value = 333
$ws4.Range("F8").Value = 329
$ws4.Range("F9").Value = 1339
$ws4.Range("F12").Value = 283
$ws4.Range("F13").Value = 496
$ws4.Range("F16").Value = 545
$ws4.Range("F17").Value = 8261
$ws4.Range("F18").Value = 359
$ws4.Range("F21").Value = 236
$ws4.Range("F25").Value = 563
$ws4.Range("F27").Value = 1147
$ws4.Range("F28").Value = 992
$ws4.Range("F29").Value = 1916
$ws4.Range("F30").Value = 1457
$ws4.Range("F31").Value = 238
$ws4.Range("F35").Value = 18
$ws4.Range("F38").Value = 297
$ws4.Range("F39").Value = 51
$ws4.Range("F40").Value = 213
$ws4.Range("F41").Value = 388
$ws4.Range("F46").Value = 44

$wb.Save()
